$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are plain text (e.g. "35.280.82"), but assigning such
# strings directly would make Excel auto-convert them to numbers. Force the cell
# to Text format before writing, then restore the default "Normal" style so the
# cell keeps no extra formatting (matching the original workbook).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.280.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.909.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("E5").Value = "  +8.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "255.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.65%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -1.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.371"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0759"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.14%  "

$ws.Range("E12").Value = "  -0.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.187.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.728"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.15%  "

$ws.Range("E16").Value = "  +1.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.902.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.279.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0850"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.66%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.35%  "

$ws.Range("E28").Value = "  +2.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.132"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128.93"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.81%  "

$ws.Range("E33").Value = "  +14.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +23.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0589"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.03%  "

$ws.Range("E37").Value = "  -0.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.910"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0218"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.09%  "

$ws.Range("E43").Value = "  +1.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0652"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.337.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.79%  "

$ws.Range("E47").Value = "  +1.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0749"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.07%  "
